$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.698.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4724"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3970"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08060"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.028"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.879.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.973"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.205"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06594"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.708.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.517"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.105.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.100"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09541"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9626"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.479"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.622"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.308"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06132"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02258"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.228"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6014"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1900"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.245"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.405"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.939"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06830"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
